$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old placeholder content (single A1 cell).
$ws.Cells.Clear()

# Fill the "tag"/count column (D, plus I/J) first so the shared-string
# table picks up "rr"/"ee" before the key/value columns below.
$ws.Range("D10").Value = "rr"
$ws.Range("D11").Value = "ee"
$ws.Range("D12").Value = "ee"
$ws.Range("D13").Value = "ee"
$ws.Range("I13").Value = "ee"
$ws.Range("J13").Value = "ee"
$ws.Range("I14").Value = "ee"
$ws.Range("J14").Value = "ee"

# Now fill the key/value columns (B, C).
$ws.Range("B10").Value = "test"
$ws.Range("C10").Value = "val"

$ws.Range("B11").Value = "this"
$ws.Range("C11").Value = "is"
$ws.Range("C11").NumberFormat = "d-mmm"

$ws.Range("B12").Value = "a"
$ws.Range("C12").Value = "test"
$ws.Range("C12").NumberFormat = "d-mmm"

$ws.Range("B13").Value = "tab"
$ws.Range("C13").Value = "table"

# Column A width, matching the authored layout.
$ws.Columns.Item(1).ColumnWidth = 24.83203125

# Selection / view state.
$ws.Range("B17").Select() | Out-Null
